# Apply change described by the diff:
# Insert a new row before row 204 (shifting existing rows 204-207 down to 205-208),
# and populate the new row 204 with data mirroring the old row 204 but with
# updated Fecha (D) and Volumen (J) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204; this shifts rows 204:207 down to 205:208
$ws.Rows.Item(204).Insert()

# New row 204 values (copy of old row 204's static data, with new D and J)
$ws.Cells.Item(204, 1).Value = 10
$ws.Cells.Item(204, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(204, 3).Value = "La Araucanía"
$ws.Cells.Item(204, 4).Value = 44448
$ws.Cells.Item(204, 5).Value = 9
$ws.Cells.Item(204, 6).Value = 100112040
$ws.Cells.Item(204, 7).Value = "Cilantro"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 85
$ws.Cells.Item(204, 11).Value = 4000
$ws.Cells.Item(204, 12).Value = 4000
$ws.Cells.Item(204, 13).Value = 4000
$ws.Cells.Item(204, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(204, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(204, 16).Value = 2000
$ws.Cells.Item(204, 17).Value = 2
$ws.Cells.Item(204, 18).Value = "Hortaliza"

# Apply the same style (numeric date format) to the new D204 cell as used elsewhere in column D
$ws.Cells.Item(204, 4).NumberFormat = $ws.Cells.Item(205, 4).NumberFormat
